# fix: unique command names in XLSX - prefix protocol name to each step
#
# For every "protocol" worksheet below, Column A (rows 2..last-used-row) holds a
# short step/command name ("Step4 Seed", "delay2", "h8", ...). Those names collide
# across worksheets (every *1/*2 pair repeats "Step1 ...", "Step2 ...", etc.), so
# this prepends the worksheet's own name + a space to each of those Column A
# values, making every command name unique workbook-wide. Row 1 (the
# Name/Text/Note/*Guidelines header) and columns B/C/D are left untouched, as are
# the five non-protocol sheets (MiaBrooksJourney, NRWaves, PersonalMiaBrooks,
# PositiveSpin, ReEngagement), which don't use this naming scheme.

$wb = $excel.ActiveWorkbook

$xlUp = -4162

$protocolSheets = @(
    "price1", "price2",
    "discount1", "discount2",
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol",
    "dickpic",
    "boosters"
)

foreach ($sheetName in $protocolSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $prefix = $sheetName + " "

    # Last used row in column A (mirrors Excel's Ctrl+Up from the bottom of the sheet).
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Text
        if ($current -and -not $current.StartsWith($prefix)) {
            $cell.Value = $prefix + $current
        }
    }
}
